# Update the "想去人数" (number of people interested) column (F) values
# on the "展览" and "全部类型" worksheets, reflecting refreshed scrape data.

$wb = $excel.ActiveWorkbook

# Worksheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5
$ws1.Range("F4").Value = 209
$ws1.Range("F5").Value = 2546
$ws1.Range("F6").Value = 1836
$ws1.Range("F7").Value = 353
$ws1.Range("F8").Value = 105
$ws1.Range("F9").Value = 892

# Worksheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5
$ws4.Range("F4").Value = 209
$ws4.Range("F5").Value = 2546
$ws4.Range("F6").Value = 1836
$ws4.Range("F7").Value = 353
$ws4.Range("F9").Value = 105
$ws4.Range("F10").Value = 892
